# Updates cryptocurrency Price (col D) and Volume(1h) (col E) values on
# Sheet1, matching the latest scrape results from the GitHub Actions job.
# A leading apostrophe is used for values that look like plain numbers so
# Excel keeps them as literal text (e.g. "0.9998", "244.85") instead of
# converting them to floating point numbers - this matches the source data
# which stores every Price/Volume cell as a string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.893.77'
$ws.Range("E2").Value = '  -1.18%  '
$ws.Range("D3").Value = '1.830.31'
$ws.Range("E3").Value = '  -1.64%  '
$ws.Range("D4").Value = '''0.9998'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''244.85'
$ws.Range("E5").Value = '  +1.03%  '
$ws.Range("D6").Value = '''0.6924'
$ws.Range("E6").Value = '  -0.89%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '''0.07665'
$ws.Range("E8").Value = '  -2.13%  '
$ws.Range("D9").Value = '''0.3038'
$ws.Range("E9").Value = '  -2.71%  '
$ws.Range("D10").Value = '''23.26'
$ws.Range("E10").Value = '  -3.27%  '
$ws.Range("E11").Value = '  +0.11%  '
$ws.Range("D12").Value = '''92.95'
$ws.Range("E12").Value = '  +1.06%  '
$ws.Range("D13").Value = '1.833.27'
$ws.Range("E13").Value = '  -1.59%  '
$ws.Range("D14").Value = '''5.090'
$ws.Range("E14").Value = '  -0.89%  '
$ws.Range("D15").Value = '''0.6808'
$ws.Range("E15").Value = '  -1.85%  '
$ws.Range("D16").Value = '''6.524'
$ws.Range("E16").Value = '  -1.18%  '
$ws.Range("D17").Value = '''0.000008233'
$ws.Range("E17").Value = '  -3.35%  '
$ws.Range("D18").Value = '28.917.77'
$ws.Range("E18").Value = '  -1.22%  '
$ws.Range("D19").Value = '''239.48'
$ws.Range("E19").Value = '  -3.49%  '
$ws.Range("D20").Value = '2.074.48'
$ws.Range("E20").Value = '  -2.13%  '
$ws.Range("D21").Value = '''12.65'
$ws.Range("E21").Value = '  -2.22%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").Value = '''7.446'
$ws.Range("E23").Value = '  -1.75%  '
$ws.Range("D24").Value = '''1.000'
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").Value = '''0.1498'
$ws.Range("E25").Value = '  -2.79%  '
$ws.Range("D26").Value = '''157.91'
$ws.Range("E26").Value = '  -1.83%  '
$ws.Range("D27").Value = '''8.727'
$ws.Range("E27").Value = '  -2.07%  '
$ws.Range("D28").Value = '''18.14'
$ws.Range("E28").Value = '  -2.53%  '
$ws.Range("D29").Value = '''1.539'
$ws.Range("E29").Value = '  -2.49%  '
$ws.Range("D30").Value = '''4.222'
$ws.Range("E30").Value = '  -1.47%  '
$ws.Range("D31").Value = '''4.131'
$ws.Range("E31").Value = '  -2.59%  '
$ws.Range("D33").Value = '''0.05108'
$ws.Range("E33").Value = '  -2.38%  '
$ws.Range("D34").Value = '''0.7744'
$ws.Range("E34").Value = '  +1.66%  '
$ws.Range("D35").Value = '''1.846'
$ws.Range("E35").Value = '  -1.85%  '
$ws.Range("D36").Value = '''1.139'
$ws.Range("D37").Value = '''2.693'
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("D38").Value = '1.273.71'
$ws.Range("E38").Value = '  +2.36%  '
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("D40").Value = '''2.697'
$ws.Range("E40").Value = '  -1.62%  '
$ws.Range("D41").Value = '''0.9537'
$ws.Range("E41").Value = '  +5.79%  '
$ws.Range("D42").Value = '''6.144'
$ws.Range("E42").Value = '  +4.01%  '
$ws.Range("D43").Value = '''106.74'
$ws.Range("E43").Value = '  -3.29%  '
$ws.Range("D44").Value = '''0.9997'
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").Value = '''9.668'
$ws.Range("E45").Value = '  +1.06%  '
$ws.Range("E46").Value = '  -1.47%  '
$ws.Range("D47").Value = '''0.5166'
$ws.Range("E47").Value = '  -0.32%  '
$ws.Range("D48").Value = '1.974.92'
$ws.Range("E48").Value = '  -1.90%  '
$ws.Range("D49").Value = '''63.72'
$ws.Range("E49").Value = '  -7.59%  '
$ws.Range("D50").Value = '''1.748'
$ws.Range("E50").Value = '  -1.26%  '
$ws.Range("D51").Value = '''6.958'
$ws.Range("E51").Value = '  -0.90%  '
